# Fixed naive component forecaster bug - Presentation state 11.02.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K24").Value = -10.42667766443478
$ws.Range("J25").Value = -10.42367604222312
$ws.Range("K25").Value = 3.722140930129188
$ws.Range("I26").Value = -10.52581323285086
$ws.Range("J26").Value = 3.620003739501452
$ws.Range("K26").Value = 5.647005785655746
$ws.Range("H27").Value = -10.5593423204977
$ws.Range("I27").Value = 3.586474651854613
$ws.Range("J27").Value = 5.613476698008907
$ws.Range("K27").Value = -2.976209967492681
$ws.Range("G28").Value = -10.46993154308055
$ws.Range("H28").Value = 3.675885429271756
$ws.Range("I28").Value = 5.702887475426051
$ws.Range("J28").Value = -2.886799190075538
$ws.Range("K28").Value = -5.861013661827103
$ws.Range("F29").Value = -10.51018406657204
$ws.Range("G29").Value = 3.635632905780271
$ws.Range("H29").Value = 5.662634951934566
$ws.Range("I29").Value = -2.927051713567023
$ws.Range("J29").Value = -5.901266185318589
$ws.Range("K29").Value = 2.056557596745962
$ws.Range("E30").Value = -10.51800434283648
$ws.Range("F30").Value = 3.627812629515833
$ws.Range("G30").Value = 5.654814675670127
$ws.Range("H30").Value = -2.934871989831461
$ws.Range("I30").Value = -5.909086461583027
$ws.Range("J30").Value = 2.048737320481524
$ws.Range("K30").Value = 0.4327147679773617
$ws.Range("D31").Value = -10.59054831629411
$ws.Range("E31").Value = 3.555268656058203
$ws.Range("F31").Value = 5.582270702212497
$ws.Range("G31").Value = -3.007415963289091
$ws.Range("H31").Value = -5.981630435040657
$ws.Range("I31").Value = 1.976193347023894
$ws.Range("J31").Value = 0.3601707945197319
$ws.Range("K31").Value = -0.5925567769544813
$ws.Range("C32").Value = -10.92245430020549
$ws.Range("D32").Value = 3.223362672146816
$ws.Range("E32").Value = 5.25036471830111
$ws.Range("F32").Value = -3.339321947200478
$ws.Range("G32").Value = -6.313536418952044
$ws.Range("H32").Value = 1.644287363112507
$ws.Range("I32").Value = 0.02826481060834481
$ws.Range("J32").Value = -0.9244627608658684
$ws.Range("K32").Value = -1.144358446982199
$ws.Range("B33").Value = -11.95530922846234
$ws.Range("C33").Value = 2.19050774388997
$ws.Range("D33").Value = 4.217509790044264
$ws.Range("E33").Value = -4.372176875457324
$ws.Range("F33").Value = -7.34639134720889
$ws.Range("G33").Value = 0.611432434855661
$ws.Range("H33").Value = -1.004590117648501
$ws.Range("I33").Value = -1.957317689122714
$ws.Range("J33").Value = -2.177213375239045
$ws.Range("K33").Value = -0.1353519709171249
$ws.Range("B34").Value = 3.268665508691584
$ws.Range("C34").Value = 5.295667554845878
$ws.Range("D34").Value = -3.29401911065571
$ws.Range("E34").Value = -6.268233582407277
$ws.Range("F34").Value = 1.689590199657275
$ws.Range("G34").Value = 0.07356764715311254
$ws.Range("H34").Value = -0.8791599243211007
$ws.Range("I34").Value = -1.099055610437431
$ws.Range("J34").Value = 0.9428057938844886
$ws.Range("K34").Value = 1.168734197848267
$ws.Range("B35").Value = 4.950828583147882
$ws.Range("C35").Value = -3.638858082353706
$ws.Range("D35").Value = -6.613072554105272
$ws.Range("E35").Value = 1.344751227959279
$ws.Range("F35").Value = -0.2712713245448828
$ws.Range("G35").Value = -1.223998896019096
$ws.Range("H35").Value = -1.443894582135427
$ws.Range("I35").Value = 0.5979668221864933
$ws.Range("J35").Value = 0.8238952261502714
$ws.Range("K35").Value = -0.4533206709007789
$ws.Range("B36").Value = -4.068997055914149
$ws.Range("C36").Value = -7.043211527665715
$ws.Range("D36").Value = 0.9146122543988354
$ws.Range("E36").Value = -0.7014102981053266
$ws.Range("F36").Value = -1.65413786957954
$ws.Range("G36").Value = -1.87403355569587
$ws.Range("H36").Value = 0.1678278486260495
$ws.Range("I36").Value = 0.3937562525898276
$ws.Range("J36").Value = -0.8834596444612227
$ws.Range("K36").Value = -0.5704286793306007
$ws.Range("B37").Value = -6.746843049765346
$ws.Range("C37").Value = 1.210980732299205
$ws.Range("D37").Value = -0.4050418202049567
$ws.Range("E37").Value = -1.35776939167917
$ws.Range("F37").Value = -1.577665077795501
$ws.Range("G37").Value = 0.4641963265264194
$ws.Range("H37").Value = 0.6901247304901975
$ws.Range("I37").Value = -0.5870911665608528
$ws.Range("J37").Value = -0.2740602014302308
$ws.Range("K37").Value = -0.7764753152880101
$ws.Range("B38").Value = 1.785526113499941
$ws.Range("C38").Value = 0.1695035609957792
$ws.Range("D38").Value = -0.783224010478434
$ws.Range("E38").Value = -1.003119696594765
$ws.Range("F38").Value = 1.038741707727155
$ws.Range("G38").Value = 1.264670111690933
$ws.Range("H38").Value = -0.01254578536011693
$ws.Range("I38").Value = 0.3004851797705051
$ws.Range("J38").Value = -0.2019299340872742
$ws.Range("K38").Value = 0.4497033133215503
$ws.Range("B39").Value = -0.03967578553597981
$ws.Range("C39").Value = -0.992403357010193
$ws.Range("D39").Value = -1.212299043126524
$ws.Range("E39").Value = 0.8295623611953963
$ws.Range("F39").Value = 1.055490765159175
$ws.Range("G39").Value = -0.2217251318918759
$ws.Range("H39").Value = 0.09130583323874608
$ws.Range("I39").Value = -0.4111092806190332
$ws.Range("J39").Value = 0.2405239667897913
$ws.Range("K39").Value = 0.0476243477478103
$ws.Range("B40").Value = -0.9569421940374809
$ws.Range("C40").Value = -1.176837880153812
$ws.Range("D40").Value = 0.8650235241681083
$ws.Range("E40").Value = 1.090951928131886
$ws.Range("F40").Value = -0.1862639689191639
$ws.Range("G40").Value = 0.1267669962114581
$ws.Range("H40").Value = -0.3756481176463212
$ws.Range("I40").Value = 0.2759851297625033
$ws.Range("J40").Value = 0.08308551072052231
$ws.Range("K40").Value = -0.4276559378643726
$ws.Range("B41").Value = -1.000053692027457
$ws.Range("C41").Value = 1.041807712294463
$ws.Range("D41").Value = 1.267736116258241
$ws.Range("E41").Value = -0.009479780792808909
$ws.Range("F41").Value = 0.3035511843378131
$ws.Range("G41").Value = -0.1988639295199662
$ws.Range("H41").Value = 0.4527693178888583
$ws.Range("I41").Value = 0.2598696988468773
$ws.Range("J41").Value = -0.2508717497380176
$ws.Range("K41").Value = -0.04547101714081951
$ws.Range("B42").Value = 1.38437779088911
$ws.Range("C42").Value = 1.610306194852888
$ws.Range("D42").Value = 0.3330902978018372
$ws.Range("E42").Value = 0.6461212629324592
$ws.Range("F42").Value = 0.1437061490746799
$ws.Range("G42").Value = 0.7953393964835044
$ws.Range("H42").Value = 0.6024397774415234
$ws.Range("I42").Value = 0.09169832885662849
$ws.Range("J42").Value = 0.2970990614538266
$ws.Range("K42").Value = -0.1093868805094843
$ws.Range("B43").Value = 2.568311565460754
$ws.Range("C43").Value = 1.291095668409704
$ws.Range("D43").Value = 1.604126633540326
$ws.Range("E43").Value = 1.101711519682546
$ws.Range("F43").Value = 1.753344767091371
$ws.Range("G43").Value = 1.56044514804939
$ws.Range("H43").Value = 1.049703699464495
$ws.Range("I43").Value = 1.255104432061693
$ws.Range("J43").Value = 0.8486184900983822
$ws.Range("K43").Value = 1.283722647366724
$ws.Range("B44").Value = 0.06590207884413968
$ws.Range("C44").Value = 0.3789330439747617
$ws.Range("D44").Value = -0.1234820698830176
$ws.Range("E44").Value = 0.5281511775258069
$ws.Range("F44").Value = 0.3352515584838259
$ws.Range("G44").Value = -0.175489890101069
$ws.Range("H44").Value = 0.02991084249612908
$ws.Range("I44").Value = -0.3765750994671818
$ws.Range("J44").Value = 0.05852905780115958
$ws.Range("B45").Value = 0.7005922962224502
$ws.Range("C45").Value = 0.198177182364671
$ws.Range("D45").Value = 0.8498104297734954
$ws.Range("E45").Value = 0.6569108107315145
$ws.Range("F45").Value = 0.1461693621466196
$ws.Range("G45").Value = 0.3515700947438177
$ws.Range("H45").Value = -0.05491584721949322
$ws.Range("I45").Value = 0.3801883100488482
$ws.Range("B46").Value = 0.6769520741007256
$ws.Range("C46").Value = 1.32858532150955
$ws.Range("D46").Value = 1.135685702467569
$ws.Range("E46").Value = 0.6249442538826742
$ws.Range("F46").Value = 0.8303449864798722
$ws.Range("G46").Value = 0.4238590445165613
$ws.Range("H46").Value = 0.8589632017849027
$ws.Range("B47").Value = 1.004139131231226
$ws.Range("C47").Value = 0.8112395121892447
$ws.Range("D47").Value = 0.3004980636043498
$ws.Range("E47").Value = 0.5058987962015479
$ws.Range("F47").Value = 0.09941285423823698
$ws.Range("G47").Value = 0.5345170115065784
$ws.Range("B48").Value = 0.1605125027665679
$ws.Range("C48").Value = -0.350228945818327
$ws.Range("D48").Value = -0.1448282132211289
$ws.Range("E48").Value = -0.5513141551844398
$ws.Range("F48").Value = -0.1162099979160984
$ws.Range("B49").Value = -0.1674057910879981
$ws.Range("C49").Value = 0.03799494150919998
$ws.Range("D49").Value = -0.3684910004541109
$ws.Range("E49").Value = 0.06661315681423048
$ws.Range("B50").Value = -0.120096269830519
$ws.Range("C50").Value = -0.5265822117938299
$ws.Range("D50").Value = -0.09147805452548849
$ws.Range("B51").Value = -0.4929177499008816
$ws.Range("C51").Value = -0.05781359263254018
$ws.Range("B52").Value = -0.1189820961875583
